$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its exact literal text representation
# (e.g. "1.000", "0.07980") instead of being auto-converted to a number
# by Excel when the new values are assigned below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '23.187.55'
$ws.Range('E2').Value = '  -3.28%  '

$ws.Range('D3').Value = '1.604.56'
$ws.Range('E3').Value = '  -2.99%  '

$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('E5').Value = '  +0.04%  '

$ws.Range('D6').Value = '302.79'
$ws.Range('E6').Value = '  -2.31%  '

$ws.Range('D7').Value = '0.3789'
$ws.Range('E7').Value = '  -2.47%  '

$ws.Range('D8').Value = '0.3670'
$ws.Range('E8').Value = '  -4.01%  '

$ws.Range('D9').Value = '50.10'
$ws.Range('E9').Value = '  -4.73%  '

$ws.Range('D10').Value = '1.271'
$ws.Range('E10').Value = '  -6.02%  '

$ws.Range('D11').Value = '0.08173'
$ws.Range('E11').Value = '  -3.36%  '

$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  +0.10%  '

$ws.Range('D13').Value = '23.21'
$ws.Range('E13').Value = '  -2.80%  '

$ws.Range('D14').Value = '6.633'
$ws.Range('E14').Value = '  -6.48%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '7.453'
$ws.Range('E15').Value = '  -6.44%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.00001261'
$ws.Range('E16').Value = '  -3.73%  '

$ws.Range('D17').Value = '1.610.79'
$ws.Range('E17').Value = '  -3.49%  '

$ws.Range('E18').Value = '  -3.20%  '

$ws.Range('D19').Value = '0.06840'
$ws.Range('E19').Value = '  -2.11%  '

$ws.Range('D20').Value = '18.35'
$ws.Range('E20').Value = '  -7.04%  '

$ws.Range('D21').Value = '6.592'
$ws.Range('E21').Value = '  -5.49%  '

$ws.Range('E22').Value = '  -0.06%  '

$ws.Range('D23').Value = '13.09'
$ws.Range('E23').Value = '  -4.89%  '

$ws.Range('D24').Value = '23.187.05'
$ws.Range('E24').Value = '  -3.29%  '

$ws.Range('D25').Value = '2.350'
$ws.Range('E25').Value = '  -4.40%  '

$ws.Range('D26').Value = '2.821'
$ws.Range('E26').Value = '  -5.17%  '

$ws.Range('D27').Value = '21.15'
$ws.Range('E27').Value = '  -4.50%  '

$ws.Range('D28').Value = '151.60'
$ws.Range('E28').Value = '  -0.56%  '

$ws.Range('D29').Value = '5.302'
$ws.Range('E29').Value = '  -2.52%  '

$ws.Range('D30').Value = '132.96'
$ws.Range('E30').Value = '  -4.15%  '

$ws.Range('D31').Value = '2.455'
$ws.Range('E31').Value = '  -2.91%  '

$ws.Range('D32').Value = '6.888'
$ws.Range('E32').Value = '  -13.01%  '

$ws.Range('D33').Value = '1.789.46'
$ws.Range('E33').Value = '  -3.17%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.07772'
$ws.Range('E34').Value = '  -4.19%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').Value = '0.9545'
$ws.Range('E35').Value = '  -8.26%  '

$ws.Range('D36').Value = '0.02767'
$ws.Range('E36').Value = '  -6.07%  '

$ws.Range('D37').Value = '6.292'
$ws.Range('E37').Value = '  -6.98%  '

$ws.Range('D38').Value = '0.2562'
$ws.Range('E38').Value = '  -4.27%  '

$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = '0.08957'
$ws.Range('E39').Value = '  -1.88%  '

$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = '10.18'
$ws.Range('E40').Value = '  -5.53%  '

$ws.Range('D41').Value = '1.393'

$ws.Range('D42').Value = '0.7170'
$ws.Range('E42').Value = '  -5.67%  '

$ws.Range('D43').Value = '12.75'
$ws.Range('E43').Value = '  -5.47%  '

$ws.Range('D44').Value = '15.68'
$ws.Range('E44').Value = '  -4.79%  '

$ws.Range('D45').Value = '0.6693'
$ws.Range('E45').Value = '  -4.12%  '

$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').Value = '2.324'
$ws.Range('E46').Value = '  -5.88%  '

$ws.Range('B47').Value = 'Frax'
$ws.Range('C47').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D47').Value = '1.000'
$ws.Range('E47').Value = '  +0.00%  '

$ws.Range('D48').Value = '4.000'
$ws.Range('E48').Value = '  -2.31%  '

$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '133.24'
$ws.Range('E49').Value = '  -1.32%  '

$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').Value = '0.07980'
$ws.Range('E50').Value = '  -4.13%  '

$ws.Range('D51').Value = '1.219'
$ws.Range('E51').Value = '  +1.41%  '

